$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '34.446.11'
$ws.Range("E2").Value = '  +0.17%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.799.29'
$ws.Range("E3").Value = '  -0.11%  '
$ws.Range("E4").Value = '  +0.58%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '224.47'
$ws.Range("E5").Value = '  -0.94%  '
$ws.Range("E6").Value = '  +0.75%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.01'
$ws.Range("E7").Value = '  +0.56%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '40.24'
$ws.Range("E8").Value = '  +10.96%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.289'
$ws.Range("E9").Value = '  -1.22%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0666'
$ws.Range("E10").Value = '  -3.22%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0993'
$ws.Range("E11").Value = '  +3.00%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '2.063.80'
$ws.Range("E12").Value = '  +0.13%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.814.29'
$ws.Range("E13").Value = '  +0.63%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '10.79'
$ws.Range("E14").Value = '  -4.25%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '34.490.66'
$ws.Range("E15").Value = '  +0.36%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.628'
$ws.Range("E16").Value = '  -1.21%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '4.37'
$ws.Range("E17").Value = '  -0.86%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '67.27'
$ws.Range("E18").Value = '  -3.37%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '239.34'
$ws.Range("E19").Value = '  -1.46%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0764'
$ws.Range("E20").Value = '  -2.29%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.01'
$ws.Range("E21").Value = '  -3.06%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.01'
$ws.Range("E22").Value = '  +0.33%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.06'
$ws.Range("E23").Value = '  -2.27%  '
$ws.Range("E24").Value = '  -4.14%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '172.10'
$ws.Range("E25").Value = '  +0.76%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.62'
$ws.Range("E26").Value = '  -4.61%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '17.25'
$ws.Range("E27").Value = '  -0.73%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.120'
$ws.Range("E28").Value = '  -0.21%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.01'
$ws.Range("E29").Value = '  +0.63%  '
$ws.Range("E30").Value = '  -0.97%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.76'
$ws.Range("E31").Value = '  -1.36%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0513'
$ws.Range("E32").Value = '  -0.85%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.82'
$ws.Range("E33").Value = '  -2.61%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.77'
$ws.Range("E34").Value = '  -0.89%  '
$ws.Range("B35").Value = 'Maker'
$ws.Range("C35").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.308.44'
$ws.Range("E35").Value = '  -4.61%  '
$ws.Range("B36").Value = 'ImmutableX'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.639'
$ws.Range("E36").Value = '  -1.97%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.05'
$ws.Range("E37").Value = '  +0.24%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '14.90'
$ws.Range("E38").Value = '  +11.58%  '
$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0187'
$ws.Range("E39").Value = '  +0.52%  '
$ws.Range("B40").Value = 'WEMIXToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.25'
$ws.Range("E40").Value = '  +7.17%  '
$ws.Range("B41").Value = 'RenderToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.33'
$ws.Range("E41").Value = '  -0.56%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '83.54'
$ws.Range("E42").Value = '  +2.06%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.44'
$ws.Range("E43").Value = '  +1.06%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.79'
$ws.Range("E44").Value = '  -0.11%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.934'
$ws.Range("E45").Value = '  -1.05%  '
$ws.Range("E46").Value = '  +4.97%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.967.63'
$ws.Range("E47").Value = '  +0.27%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '5.75'
$ws.Range("E48").Value = '  -1.76%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.01'
$ws.Range("E49").Value = '  +0.53%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '100.50'
$ws.Range("E50").Value = '  -1.82%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0608'
$ws.Range("E51").Value = '  +0.47%  '
